# Regenerate save_data to use K instead of Strike#, recalc std/mean, write s_vals.
# In this workbook the change surfaces as updated values in column G ("K")
# for the game log rows (rows 2-13 on the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value (replacing old Strike# derived value)
$kValues = @{
    2  = 5
    3  = 5
    4  = 2
    5  = 2
    6  = 3
    7  = 3
    8  = 1
    9  = 1
    10 = 8
    11 = 3
    12 = 2
    13 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
